# Junction_Flooding_104.xlsx edit:
#  - round row-5 data cells (B5:AH5) to 2 decimal places ("custom accuracy")
#  - drop row 6 entirely (was a duplicate/extra reading)
#  - tighten columns J (10) and AB (28) by one character unit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Round the numeric readings in row 5 to 2 decimal places ---------
$lastCol = 34   # column AH
for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $excel.WorksheetFunction.Round($cell.Value2, 2)
    }
}

# --- 2. Remove row 6 (extra data point no longer wanted) -----------------
$ws.Rows.Item(6).Delete()

# --- 3. Narrow columns J and AB from 8 to 7 character units --------------
# COM ColumnWidth and the stored OOXML character-width differ by 5/6, so
# subtract that offset to land exactly on the target stored width of 7.
$ws.Columns.Item(10).ColumnWidth = 7 - 5/6
$ws.Columns.Item(28).ColumnWidth = 7 - 5/6
